# "Generate Report for Handoff"
# Re-running the handoff report generation advances the "Latest Handoff
# Datetime" timestamps for the files that just got handed off (status
# "Ready for handoff", rows 8-13 on each language sheet) and stamps their
# Priority column with "ht" (handoff type). The Overview sheet's
# "Latest HO Xliff Generate Date" column mirrors the de-de sheet's
# handoff datetime for those same rows.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# zh-cn: mark rows 8-13 (Ready for handoff) with priority "ht" and
# refresh their Latest Handoff Datetime.
$wsZhCn.Range("E8:E13").Value = "ht"
$wsZhCn.Range("H8:H13").Value = "2016-08-31 08:25:16"

# de-de: same update, with its own regenerated handoff datetime.
$wsDeDe.Range("E8:E13").Value = "ht"
$wsDeDe.Range("H8:H13").Value = "2016-08-31 08:25:22"

# Overview: Latest HO Xliff Generate Date follows the de-de handoff time.
$wsOverview.Range("G8:G13").Value = "2016-08-31 08:25:22"
